# "Completed some BFS practices." — add two new Binary Tree Level Order
# Traversal entries to the top of the "Breadth First Search" sheet, pushing
# the existing two rows (Populating Next Right Pointers ...) down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Breadth First Search")

# ---------------------------------------------------------------------
# 1. Make room: insert two fresh rows above the current row 2. The old
#    row 2/3 data (Populating Next Right Pointers ...) ends up at row 4/5.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# Helper: apply the plain "data cell" look (bordered, centered, 10pt
# Microsoft YaHei) used across this sheet to a range.
function Set-DataStyle($rng) {
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.Borders.LineStyle = 1
    $rng.Font.Name = "Microsoft YaHei"
    $rng.Font.Size = 10
}

# ---------------------------------------------------------------------
# 2. Row 2: "Binary Tree Level Order Traversal"
# ---------------------------------------------------------------------
Set-DataStyle($ws.Range("A2:B2"))
Set-DataStyle($ws.Range("D2:H2"))
Set-DataStyle($ws.Range("K2"))

$ws.Range("A2").Value = 102
$ws.Range("B2").Value = "Binary Tree Level Order Traversal"
$ws.Range("D2").Value = "Medium"
$ws.Range("F2").Value = "Binary Tree"
$ws.Range("G2").Value = "O(n)"
$ws.Range("H2").Value = "O(n)"

# E2: wrapped note text
$e2 = $ws.Range("E2")
Set-DataStyle($e2)
$e2.WrapText = $true
$e2.Value = "If needed, use .copy() to prevent`nstorage modification during iteration."

# C2: hyperlink-styled cell (Calibri, blue, underlined, via builtin style)
$c2 = $ws.Range("C2")
$c2.Style = "Hyperlink"
Set-DataStyle($c2)
$c2.Font.Name = "Calibri"
$c2.Font.Size = 11
$c2.Value = "https://leetcode.com/problems/binary-tree-level-order-traversal/description/"

# I2/J2: date-formatted cells
$ij2 = $ws.Range("I2:J2")
Set-DataStyle($ij2)
$ij2.NumberFormat = "yyyy""年""m""月""d""日"";@"
$ws.Range("I2").Value = 45489

# ---------------------------------------------------------------------
# 3. Row 3: "	Binary Tree Level Order Traversal II"
# ---------------------------------------------------------------------
Set-DataStyle($ws.Range("A3"))
Set-DataStyle($ws.Range("D3:H3"))
Set-DataStyle($ws.Range("K3"))

$ws.Range("A3").Value = 107
$ws.Range("D3").Value = "Medium"
$ws.Range("F3").Value = "Binary Tree"
$ws.Range("G3").Value = "O(n)"
$ws.Range("H3").Value = "O(n)"

# B3: wrapped title (longer name, gets the wrap-text look like other
# tab-prefixed "II" entries in this workbook)
$b3 = $ws.Range("B3")
Set-DataStyle($b3)
$b3.WrapText = $true
$b3.Value = "`tBinary Tree Level Order Traversal II"

# E3: wrapped note text (same note as row 2)
$e3 = $ws.Range("E3")
Set-DataStyle($e3)
$e3.WrapText = $true
$e3.Value = "If needed, use .copy() to prevent`nstorage modification during iteration."

# C3: hyperlink-styled cell (Calibri, blue, underlined, via builtin style)
$c3 = $ws.Range("C3")
$c3.Style = "Hyperlink"
Set-DataStyle($c3)
$c3.Font.Name = "Calibri"
$c3.Font.Size = 11
$c3.Value = "https://leetcode.com/problems/binary-tree-level-order-traversal-ii/description/"

# I3/J3: date-formatted cells
$ij3 = $ws.Range("I3:J3")
Set-DataStyle($ij3)
$ij3.NumberFormat = "yyyy""年""m""月""d""日"";@"
$ws.Range("I3").Value = 45489

# ---------------------------------------------------------------------
# 4. Conditional formatting on the Difficulty column now spans D2:D5.
# ---------------------------------------------------------------------
foreach ($fc in $ws.Range("D2:D3").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("D2:D5"))
}

# ---------------------------------------------------------------------
# 5. Hyperlinks: rebuild them all so the two surviving links land on
#    their shifted cells (C4, C5) and the two new links land on C2/C3,
#    matching the rId1..rId4 order in the target file.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), "https://leetcode.com/problems/populating-next-right-pointers-in-each-node/")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://leetcode.com/problems/populating-next-right-pointers-in-each-node-ii/description/")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://leetcode.com/problems/binary-tree-level-order-traversal/description/")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://leetcode.com/problems/binary-tree-level-order-traversal-ii/description/")
